# Update "to go" (F column) counters across several sheets, and insert a new
# performance event row in the "演出" (Performance) sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibition)
$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performance)
$ws3 = $wb.Worksheets.Item(3)   # 本地生活 (Local life)
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All types)

# ---------------------------------------------------------------------------
# Sheet 1 (展览): bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$ws1.Cells.Item(7, 6).Value = 24
$ws1.Cells.Item(8, 6).Value = 45
$ws1.Cells.Item(9, 6).Value = 751
$ws1.Cells.Item(10, 6).Value = 2725
$ws1.Cells.Item(11, 6).Value = 2725
$ws1.Cells.Item(12, 6).Value = 21
$ws1.Cells.Item(13, 6).Value = 1835
$ws1.Cells.Item(14, 6).Value = 623
$ws1.Cells.Item(15, 6).Value = 310
$ws1.Cells.Item(17, 6).Value = 10
$ws1.Cells.Item(18, 6).Value = 6311
$ws1.Cells.Item(19, 6).Value = 244
$ws1.Cells.Item(20, 6).Value = 93
$ws1.Cells.Item(21, 6).Value = 693
$ws1.Cells.Item(27, 6).Value = 2479
$ws1.Cells.Item(33, 6).Value = 1323
$ws1.Cells.Item(39, 6).Value = 1519
$ws1.Cells.Item(40, 6).Value = 34
$ws1.Cells.Item(41, 6).Value = 1471

# ---------------------------------------------------------------------------
# Sheet 2 (演出): bump F11, then insert a brand-new event as row 15
# (shifting the former rows 15-25 down to rows 16-26).
# ---------------------------------------------------------------------------
$ws2.Cells.Item(11, 6).Value = 169

$ws2.Rows.Item(15).Insert(-4121)   # -4121 = xlShiftDown

$ws2.Cells.Item(15, 1).Value = 14
$ws2.Cells.Item(15, 2).Value = "2024-11-22"
$ws2.Cells.Item(15, 3).Value = "北京·Peder Elias 巡演"
$ws2.Cells.Item(15, 4).Value = "建国门外郎家园10号61幢一层A3-06、二层A3-06号 EAST LIVE(东郎展演中心)"
$ws2.Cells.Item(15, 5).Value = "2024.11.22 19:30-11.22 21:00"
$ws2.Cells.Item(15, 6).Value = 0
$ws2.Cells.Item(15, 7).Value = 280
$ws2.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93948"
$ws2.Cells.Item(15, 9).Value = "//i2.hdslb.com/bfs/openplatform/202410/2ElTnEKU1729737260298.jpeg"

# The "to go" counter of one of the shifted events (now row 19) also grew.
$ws2.Cells.Item(19, 6).Value = 270

# Column A is a plain running index (row number - 1); re-number it for every
# row from the insertion point through the newly extended last row.
for ($r = 15; $r -le 26; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet 3 (本地生活): bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$ws3.Cells.Item(3, 6).Value = 937
$ws3.Cells.Item(7, 6).Value = 77
$ws3.Cells.Item(8, 6).Value = 14

# ---------------------------------------------------------------------------
# Sheet 4 (全部类型): bump "want-to-go" counts (column F)
# ---------------------------------------------------------------------------
$ws4.Cells.Item(6, 6).Value = 937
$ws4.Cells.Item(11, 6).Value = 77
$ws4.Cells.Item(12, 6).Value = 77
$ws4.Cells.Item(18, 6).Value = 45
$ws4.Cells.Item(19, 6).Value = 2725
$ws4.Cells.Item(20, 6).Value = 14
$ws4.Cells.Item(22, 6).Value = 21
$ws4.Cells.Item(23, 6).Value = 169
$ws4.Cells.Item(24, 6).Value = 623
$ws4.Cells.Item(25, 6).Value = 310
$ws4.Cells.Item(27, 6).Value = 6311
$ws4.Cells.Item(28, 6).Value = 244
$ws4.Cells.Item(29, 6).Value = 93
$ws4.Cells.Item(30, 6).Value = 693
$ws4.Cells.Item(34, 6).Value = 2479
$ws4.Cells.Item(37, 6).Value = 1323
$ws4.Cells.Item(41, 6).Value = 270
$ws4.Cells.Item(48, 6).Value = 1519
$ws4.Cells.Item(49, 6).Value = 34
